$wb = $excel.ActiveWorkbook

# --- "LIST" sheet: restore column A (rows 2-10) from the backup values
#     that live in "Feuil1" column B (rows 11-19), which removes the need
#     for the orphaned "RO.ORG.001.REC" string (old row 3) entirely.
$listSheet   = $wb.Worksheets.Item("LIST")
$backupSheet = $wb.Worksheets.Item("Feuil1")

for ($r = 11; $r -le 19; $r++) {
    $v = $backupSheet.Cells.Item($r, 2).Value2
    $destCell = $listSheet.Cells.Item($r - 9, 1)
    # Row 3 (old A3) had no style; give it the same "@" text format the
    # other restored rows already carry so it matches the rest of the
    # column (cellXf s="2").
    if (($r - 9) -eq 3) {
        $destCell.NumberFormat = "@"
    }
    $destCell.Value = $v
}

# --- "Feuil1" sheet: update its stored selection to B11:B19
$backupSheet.Range("B11:B19").Select()

# --- Re-activate "LIST" (it stays the selected/visible tab) and restore
#     its selection to C8
$listSheet.Activate()
$listSheet.Range("C8").Select()
